$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New daily rows 302-328 (dates 44376-44402), matching the source diff.
$data = @(
    @(44376, 0, 1, 10.13993104846887),
    @(44377, 0, 1, 10.13993104846887),
    @(44378, 0, 1, 10.13993104846887),
    @(44379, 0, 1, 10.13993104846887),
    @(44380, 0, 0, 0),
    @(44381, 0, 0, 0),
    @(44382, 0, 0, 0),
    @(44383, 0, 0, 0),
    @(44384, 0, 0, 0),
    @(44385, 0, 0, 0),
    @(44386, 0, 0, 0),
    @(44387, 0, 0, 0),
    @(44388, 0, 0, 0),
    @(44389, 0, 0, 0),
    @(44390, 0, 0, 0),
    @(44391, 0, 0, 0),
    @(44392, 0, 0, 0),
    @(44393, 0, 0, 0),
    @(44394, 1, 1, 10.13993104846887),
    @(44395, 0, 1, 10.13993104846887),
    @(44396, 5, 6, 60.83958629081322),
    @(44397, 0, 6, 60.83958629081322),
    @(44398, 0, 6, 60.83958629081322),
    @(44399, 5, 11, 111.5392415331576),
    @(44400, 0, 11, 111.5392415331576),
    @(44401, 4, 14, 141.9590346785642),
    @(44402, 6, 20, 202.7986209693774)
)

$startRow = 302
$endRow = 328

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}

# Column A uses the same date style (s="2", center/top aligned, bordered,
# custom date number format) as the rest of the column - copy it down from
# the last existing row rather than re-building the style from scratch.
$ws.Cells.Item($startRow - 1, 1).Copy() | Out-Null
$dateRange = "A" + $startRow + ":A" + $endRow
$ws.Range($dateRange).PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
